# League bases update (10-06-2024 07:08): the match rows for each pairing
# below had their data rows swapped between the two fixtures that share the
# same match date (columns B..AD); the running index in column A is kept
# untouched for every row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 18 and row 19 (all columns except A/id)
$ws.Range("B18").Value = 6228596
$ws.Range("E18").Value = 'Young Lions'
$ws.Range("F18").Value = 'Albirex Niigata Singapore'
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 4
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 3
$ws.Range("K18").Value = 'A'
$ws.Range("L18").Value = 21
$ws.Range("M18").Value = 11
$ws.Range("N18").Value = 1.062
$ws.Range("O18").Value = 21
$ws.Range("P18").Value = 10
$ws.Range("Q18").Value = 1.083
$ws.Range("R18").Value = 3
$ws.Range("S18").Value = 1.825
$ws.Range("T18").Value = 2.025
$ws.Range("U18").Value = 4
$ws.Range("V18").Value = 1.75
$ws.Range("W18").Value = 2.125
$ws.Range("X18").Value = -1
$ws.Range("Y18").Value = -1
$ws.Range("Z18").Value = 0.08299999999999996
$ws.Range("AA18").Value = 0.825
$ws.Range("AB18").Value = -1
$ws.Range("AC18").Value = 0.75
$ws.Range("AD18").Value = -1

$ws.Range("B19").Value = 6228030
$ws.Range("E19").Value = 'Tanjong Pagar United'
$ws.Range("F19").Value = 'Lion City Sailors FC'
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = 7
$ws.Range("I19").Value = 1
$ws.Range("J19").Value = 2
$ws.Range("K19").Value = 'A'
$ws.Range("L19").Value = 7
$ws.Range("M19").Value = 6
$ws.Range("N19").Value = 1.25
$ws.Range("O19").Value = 7
$ws.Range("P19").Value = 6.5
$ws.Range("Q19").Value = 1.222
$ws.Range("R19").Value = 2
$ws.Range("S19").Value = 1.875
$ws.Range("T19").Value = 1.975
$ws.Range("U19").Value = 4.25
$ws.Range("V19").Value = 2
$ws.Range("W19").Value = 1.85
$ws.Range("X19").Value = -1
$ws.Range("Y19").Value = -1
$ws.Range("Z19").Value = 0.222
$ws.Range("AA19").Value = -1
$ws.Range("AB19").Value = 0.9750000000000001
$ws.Range("AC19").Value = 1
$ws.Range("AD19").Value = -1

# Swap row 26 and row 27 (all columns except A/id)
$ws.Range("B26").Value = 6228032
$ws.Range("E26").Value = 'Lion City Sailors FC'
$ws.Range("F26").Value = 'DPMM FC'
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 3
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 2
$ws.Range("K26").Value = 'A'
$ws.Range("L26").Value = 1.2
$ws.Range("M26").Value = 6
$ws.Range("N26").Value = 9
$ws.Range("O26").Value = 1.142
$ws.Range("P26").Value = 7.5
$ws.Range("Q26").Value = 13
$ws.Range("R26").Value = -2.25
$ws.Range("S26").Value = 1.825
$ws.Range("T26").Value = 2.025
$ws.Range("U26").Value = 4
$ws.Range("V26").Value = 1.825
$ws.Range("W26").Value = 2.025
$ws.Range("X26").Value = -1
$ws.Range("Y26").Value = -1
$ws.Range("Z26").Value = 12
$ws.Range("AA26").Value = -1
$ws.Range("AB26").Value = 1.025
$ws.Range("AC26").Value = 0
$ws.Range("AD26").Value = 0

$ws.Range("B27").Value = 6228602
$ws.Range("E27").Value = 'Tampines Rovers FC'
$ws.Range("F27").Value = 'Tanjong Pagar United'
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 1
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 'H'
$ws.Range("L27").Value = 1.2
$ws.Range("M27").Value = 6
$ws.Range("N27").Value = 9
$ws.Range("O27").Value = 1.125
$ws.Range("P27").Value = 7.5
$ws.Range("Q27").Value = 17
$ws.Range("R27").Value = -2.5
$ws.Range("S27").Value = 1.9
$ws.Range("T27").Value = 1.95
$ws.Range("U27").Value = 4.25
$ws.Range("V27").Value = 1.975
$ws.Range("W27").Value = 1.875
$ws.Range("X27").Value = 0.125
$ws.Range("Y27").Value = -1
$ws.Range("Z27").Value = -1
$ws.Range("AA27").Value = -1
$ws.Range("AB27").Value = 0.95
$ws.Range("AC27").Value = -1
$ws.Range("AD27").Value = 0.875

# Swap row 38 and row 39 (all columns except A/id)
$ws.Range("B38").Value = 6228613
$ws.Range("E38").Value = 'DPMM FC'
$ws.Range("F38").Value = 'Geylang International'
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = 2
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 1
$ws.Range("K38").Value = 'A'
$ws.Range("L38").Value = 2.5
$ws.Range("M38").Value = 3.75
$ws.Range("N38").Value = 2.25
$ws.Range("O38").Value = 2.45
$ws.Range("P38").Value = 3.6
$ws.Range("Q38").Value = 2.3
$ws.Range("R38").Value = 0
$ws.Range("S38").Value = 2
$ws.Range("T38").Value = 1.85
$ws.Range("U38").Value = 3.75
$ws.Range("V38").Value = 1.925
$ws.Range("W38").Value = 1.925
$ws.Range("X38").Value = -1
$ws.Range("Y38").Value = -1
$ws.Range("Z38").Value = 1.3
$ws.Range("AA38").Value = -1
$ws.Range("AB38").Value = 0.8500000000000001
$ws.Range("AC38").Value = -1
$ws.Range("AD38").Value = 0.925

$ws.Range("B39").Value = 6228611
$ws.Range("E39").Value = 'Albirex Niigata Singapore'
$ws.Range("F39").Value = 'Tampines Rovers FC'
$ws.Range("G39").Value = 6
$ws.Range("H39").Value = 3
$ws.Range("I39").Value = 4
$ws.Range("J39").Value = 2
$ws.Range("K39").Value = 'H'
$ws.Range("L39").Value = 1.7
$ws.Range("M39").Value = 4
$ws.Range("N39").Value = 3.6
$ws.Range("O39").Value = 1.5
$ws.Range("P39").Value = 3.8
$ws.Range("Q39").Value = 5.25
$ws.Range("R39").Value = -1.25
$ws.Range("S39").Value = 1.975
$ws.Range("T39").Value = 1.875
$ws.Range("U39").Value = 4
$ws.Range("V39").Value = 2.025
$ws.Range("W39").Value = 1.825
$ws.Range("X39").Value = 0.5
$ws.Range("Y39").Value = -1
$ws.Range("Z39").Value = -1
$ws.Range("AA39").Value = 0.9750000000000001
$ws.Range("AB39").Value = -1
$ws.Range("AC39").Value = 1.025
$ws.Range("AD39").Value = -1

# Swap row 54 and row 55 (all columns except A/id)
$ws.Range("B54").Value = 7098763
$ws.Range("E54").Value = 'Balestier Khalsa FC'
$ws.Range("F54").Value = 'Tampines Rovers FC'
$ws.Range("G54").Value = 1
$ws.Range("H54").Value = 3
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 'A'
$ws.Range("L54").Value = 5.25
$ws.Range("M54").Value = 4.2
$ws.Range("N54").Value = 1.5
$ws.Range("O54").Value = 5
$ws.Range("P54").Value = 4.5
$ws.Range("Q54").Value = 1.45
$ws.Range("R54").Value = 1.25
$ws.Range("S54").Value = 2
$ws.Range("T54").Value = 1.85
$ws.Range("U54").Value = 5
$ws.Range("V54").Value = 1.925
$ws.Range("W54").Value = 1.925
$ws.Range("X54").Value = -1
$ws.Range("Y54").Value = -1
$ws.Range("Z54").Value = 0.45
$ws.Range("AA54").Value = -1
$ws.Range("AB54").Value = 0.8500000000000001
$ws.Range("AC54").Value = -1
$ws.Range("AD54").Value = 0.925

$ws.Range("B55").Value = 7094656
$ws.Range("E55").Value = 'Tanjong Pagar United'
$ws.Range("F55").Value = 'DPMM FC'
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = 1
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 1
$ws.Range("K55").Value = 'D'
$ws.Range("L55").Value = 2.15
$ws.Range("M55").Value = 3.75
$ws.Range("N55").Value = 2.7
$ws.Range("O55").Value = 2.1
$ws.Range("P55").Value = 4.2
$ws.Range("Q55").Value = 2.625
$ws.Range("R55").Value = -0.25
$ws.Range("S55").Value = 1.925
$ws.Range("T55").Value = 1.925
$ws.Range("U55").Value = 4.25
$ws.Range("V55").Value = 1.9
$ws.Range("W55").Value = 1.95
$ws.Range("X55").Value = -1
$ws.Range("Y55").Value = 3.2
$ws.Range("Z55").Value = -1
$ws.Range("AA55").Value = -0.5
$ws.Range("AB55").Value = 0.4625
$ws.Range("AC55").Value = -1
$ws.Range("AD55").Value = 0.95
